# Auto stash before merge of "develop" and "origin/develop"
# Update shareweight values on the ETS sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Row 2: hard coal -> 0.1 for all years (columns B:AF)
$ws.Range("B2:AF2").Value = 0.1

# Row 6: onshore wind -> 0.7 for all years (columns B:AF)
$ws.Range("B6:AF6").Value = 0.7

# Row 13: lignite -> 0.1 for all years (columns B:AF)
$ws.Range("B13:AF13").Value = 0.1

# Row 14: offshore wind -> 0.3 for all years (columns B:AF)
$ws.Range("B14:AF14").Value = 0.3

# Update the selected/active view on the sheet (row 13 selected)
$ws.Activate()
$ws.Range("B13:AF13").Select()
